# Noted health charts to move
# Insert a new row at row 82 (shifting the existing rows 82-86 down to
# 83-87), mirroring the user selecting the whole row and choosing
# Insert. Row 81 ("HealthEducation/stats") then gets two new notes
# added in columns D and F about health charts ("trolleys" / patients
# on trolleys) that are to be moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the row first (matches the real editing session's selection
# before inserting), then insert a blank row above it.
$ws.Rows.Item(82).Select() | Out-Null
$ws.Rows.Item(82).Insert()

# New note cells on row 81.
$ws.Range("D81").Value = "trolleys, healthlevels"
$ws.Range("F81").Value = "Patients on Trolleys, Health Levels"

# Leave selection on the newly inserted (now blank) row, matching the
# saved view state.
$ws.Rows.Item(82).Select() | Out-Null
